# Applies the cryptos-list refresh described in the commit diff.
# All touched cells keep their original "Text" storage (matches the
# source inlineStr cells, e.g. "187.20" must not collapse to 187.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '74.172.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.627.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +6.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '187.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +13.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '582.66'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.01%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.531'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.198'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +13.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.623.32'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.56%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.69'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '73.945.64'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000189'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.22%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +6.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.52'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +11.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.611.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +6.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.18'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +29.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.80'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +9.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '366.06'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +7.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +14.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.08'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.01%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.16'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.87'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.15%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.13'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +7.10%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.32'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +9.16%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.757.46'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.30%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.01%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0942'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +11.72%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '526.14'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +17.60%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.39'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +12.35%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.66'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.13%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.74'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.72%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.26'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.118'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +7.70%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.11'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.37%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.25'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.91'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +10.23%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.67'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +8.39%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'PolygonEcosystemToken'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.325'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +7.12%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '162.89'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +23.54%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.38'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +11.37%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.18'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +7.65%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '38.98'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.82%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0846'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +17.16%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.60'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.41%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.523'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.80%  '
